$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value, per the Jan 30 2023 symbol-list refresh.
$updates = @(
    @{ Cell = "D2"; Value = "309.42" },
    @{ Cell = "E2"; Value = "-2.90%" },
    @{ Cell = "D3"; Value = "37.23" },
    @{ Cell = "E3"; Value = "-6.33%" },
    @{ Cell = "D4"; Value = "5.123" },
    @{ Cell = "E4"; Value = "-0.31%" },
    @{ Cell = "D5"; Value = "0.07836" },
    @{ Cell = "E5"; Value = "-4.74%" },
    @{ Cell = "D6"; Value = "1.961" },
    @{ Cell = "E6"; Value = "-5.00%" },
    @{ Cell = "D7"; Value = "4.395" },
    @{ Cell = "E7"; Value = "2.27%" },
    @{ Cell = "D8"; Value = "8.277" },
    @{ Cell = "E8"; Value = "-0.42%" },
    @{ Cell = "E9"; Value = "-8.68%" },
    @{ Cell = "D10"; Value = "0.9244" },
    @{ Cell = "E10"; Value = "-1.15%" },
    @{ Cell = "D11"; Value = "0.1315" },
    @{ Cell = "E11"; Value = "-3.04%" },
    @{ Cell = "D12"; Value = "0.1948" },
    @{ Cell = "E12"; Value = "-1.83%" },
    @{ Cell = "D13"; Value = "0.08965" },
    @{ Cell = "E13"; Value = "-1.21%" },
    @{ Cell = "D14"; Value = "0.03448" },
    @{ Cell = "E14"; Value = "-1.64%" },
    @{ Cell = "D15"; Value = "0.09705" },
    @{ Cell = "E15"; Value = "-1.28%" },
    @{ Cell = "D16"; Value = "0.001386" },
    @{ Cell = "E16"; Value = "-0.83%" },
    @{ Cell = "D17"; Value = "0.006040" },
    @{ Cell = "E17"; Value = "-3.98%" },
    @{ Cell = "D18"; Value = "3.587" },
    @{ Cell = "E18"; Value = "-2.60%" },
    @{ Cell = "D19"; Value = "0.3399" },
    @{ Cell = "E19"; Value = "-2.15%" },
    @{ Cell = "D20"; Value = "0.1296" },
    @{ Cell = "E20"; Value = "0.13%" },
    @{ Cell = "D21"; Value = "5.003" },
    @{ Cell = "E21"; Value = "2.09%" },
    @{ Cell = "D22"; Value = "0.2493" },
    @{ Cell = "E22"; Value = "1.90%" },
    @{ Cell = "D23"; Value = "0.02107" },
    @{ Cell = "E23"; Value = "5,176.73%" },
    @{ Cell = "D24"; Value = "0.04344" },
    @{ Cell = "E24"; Value = "0.59%" },
    @{ Cell = "D25"; Value = "0.001218" },
    @{ Cell = "E25"; Value = "-0.67%" },
    @{ Cell = "D26"; Value = "0.004517" },
    @{ Cell = "E26"; Value = "-5.33%" },
    @{ Cell = "D27"; Value = "0.0001350" },
    @{ Cell = "E27"; Value = "3.96%" },
    @{ Cell = "D39"; Value = "0.02269" },
    @{ Cell = "E39"; Value = "2.34%" },
    @{ Cell = "D40"; Value = "0.05031" },
    @{ Cell = "E40"; Value = "-3.78%" },
    @{ Cell = "D41"; Value = "0.007639" },
    @{ Cell = "E41"; Value = "-0.42%" },
    @{ Cell = "D42"; Value = "0.009802" },
    @{ Cell = "E42"; Value = "0.53%" },
    @{ Cell = "D43"; Value = "0.1352" },
    @{ Cell = "E43"; Value = "-2.72%" },
    @{ Cell = "D44"; Value = "0.001995" },
    @{ Cell = "E44"; Value = "-4.60%" },
    @{ Cell = "D45"; Value = "0.008452" },
    @{ Cell = "E45"; Value = "-8.15%" },
    @{ Cell = "D46"; Value = "0.00006775" },
    @{ Cell = "E46"; Value = "3.46%" },
    @{ Cell = "D47"; Value = "0.00000000752" },
    @{ Cell = "E47"; Value = "0.39%" },
    @{ Cell = "D48"; Value = "0.003019" },
    @{ Cell = "D49"; Value = "0.001303" },
    @{ Cell = "E49"; Value = "-22.77%" },
    @{ Cell = "D50"; Value = "0.00002105" },
    @{ Cell = "E50"; Value = "0.39%" },
    @{ Cell = "D51"; Value = "0.0002005" },
    @{ Cell = "E51"; Value = "0.39%" }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    # Force text storage (matches original inlineStr cells: price/
    # volume figures are kept as formatted text, not numbers) then
    # drop back to the default style so we do not leave a stray
    # "@" number format applied to the cell.
    $c.NumberFormat = "@"
    $c.Value = $u.Value
    $c.Style = "Normal"
}
